$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format first on any Price (column D) cells whose new value
# would otherwise be auto-parsed by Excel as a Number (losing the literal
# "thousands-dot" formatting / trailing zeros), mirroring the original
# inlineStr text storage.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price / Volume(1h) values row by row.
$ws.Range("D2").Value = "27.477.43"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "1.866.38"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "311.79"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4778"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.3762"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").Value = "0.07338"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "0.9357"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("D12").Value = "0.07844"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "1.894.60"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").Value = "5.437"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "6.554"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "90.44"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "0.000008891"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "27.546.18"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "5.119"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "1.939"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "154.57"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").Value = "18.48"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "115.56"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "4.970"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "0.08901"
$ws.Range("D31").Value = "3.333"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("D33").Value = "0.7616"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "4.605"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").Value = "2.754"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "1.124"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").Value = "0.02035"
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("D38").Value = "2.996"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "0.05265"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "0.5320"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").Value = "7.082"
$ws.Range("D42").Value = "8.557"
$ws.Range("E42").Value = "  +4.18%  "
$ws.Range("D43").Value = "0.1525"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "0.4806"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").Value = "1.013"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "102.95"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "1.655"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("D49").Value = "67.39"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "0.06081"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "0.9185"
$ws.Range("E51").Value = "  +3.40%  "
